$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.804.84'

$ws.Range('D3').Value = '2.097.71'
$ws.Range('E3').Value = '  +1.79%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '''232.25'
$ws.Range('E5').Value = '  +0.16%  '

$ws.Range('D6').Value = '''0.624'
$ws.Range('E6').Value = '  +0.08%  '

$ws.Range('B7').Value = 'Solana'
$ws.Range('C7').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D7').Value = '''58.05'
$ws.Range('E7').Value = '  +1.71%  '

$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').Value = '''0.388'
$ws.Range('E9').Value = '  +1.33%  '

$ws.Range('D10').Value = '''0.0779'
$ws.Range('E10').Value = '  +2.51%  '

$ws.Range('E11').Value = '  +2.58%  '

$ws.Range('D12').Value = '2.394.49'
$ws.Range('E12').Value = '  +1.28%  '

$ws.Range('D13').Value = '''14.49'
$ws.Range('E13').Value = '  -0.56%  '

$ws.Range('E14').Value = '  +1.86%  '

$ws.Range('E15').Value = '  -1.17%  '

$ws.Range('E16').Value = '  +2.07%  '

$ws.Range('D17').Value = '2.104.69'
$ws.Range('E17').Value = '  +2.13%  '

$ws.Range('D18').Value = '37.742.17'
$ws.Range('E18').Value = '  +1.47%  '

$ws.Range('D19').Value = '''6.15'
$ws.Range('E19').Value = '  -2.76%  '

$ws.Range('D20').Value = '''70.63'
$ws.Range('E20').Value = '  +1.87%  '

$ws.Range('D21').Value = '0.0₃0822'
$ws.Range('E21').Value = '  +1.63%  '

$ws.Range('D22').Value = '''228.08'
$ws.Range('E22').Value = '  +1.04%  '

$ws.Range('E23').Value = '  -0.08%  '

$ws.Range('E24').Value = '  -1.16%  '

$ws.Range('E25').Value = '  -0.06%  '

$ws.Range('D26').Value = '''167.88'
$ws.Range('E26').Value = '  +1.07%  '

$ws.Range('D27').Value = '''0.141'
$ws.Range('E27').Value = '  +9.57%  '

$ws.Range('D28').Value = '''8.95'
$ws.Range('E28').Value = '  +2.15%  '

$ws.Range('E29').Value = '  -1.35%  '

$ws.Range('E30').Value = '  +2.21%  '

$ws.Range('E31').Value = '  +1.36%  '

$ws.Range('D32').Value = '''4.63'
$ws.Range('E32').Value = '  +3.63%  '

$ws.Range('E33').Value = '  +1.20%  '

$ws.Range('E34').Value = '  -0.49%  '

$ws.Range('D35').Value = '''2.51'
$ws.Range('E35').Value = '  +0.42%  '

$ws.Range('E36').Value = '  +5.22%  '

$ws.Range('E37').Value = '  +4.13%  '

$ws.Range('D38').Value = '''0.999'
$ws.Range('E38').Value = '  -0.20%  '

$ws.Range('D39').Value = '''5.41'
$ws.Range('E39').Value = '  -5.02%  '

$ws.Range('D40').Value = '''0.0994'
$ws.Range('E40').Value = '  +6.26%  '

$ws.Range('E41').Value = '  -0.49%  '

$ws.Range('D42').Value = '''98.03'
$ws.Range('E42').Value = '  +1.93%  '

$ws.Range('E43').Value = '  +0.89%  '

$ws.Range('D44').Value = '1.457.92'
$ws.Range('E44').Value = '  -0.70%  '

$ws.Range('E45').Value = '  -0.74%  '

$ws.Range('E46').Value = '  +3.97%  '

$ws.Range('D47').Value = '''15.66'
$ws.Range('E47').Value = '  +3.70%  '

$ws.Range('D48').Value = '''4.04'
$ws.Range('E48').Value = '  -3.12%  '

$ws.Range('D49').Value = '''7.36'
$ws.Range('E49').Value = '  +2.78%  '

$ws.Range('E50').Value = '  +2.25%  '

$ws.Range('D51').Value = '2.288.64'
$ws.Range('E51').Value = '  +1.65%  '
